$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.507.98"
$ws.Range("E2").Value = "  -4.45%  "
$ws.Range("D3").Value = "2.914.94"
$ws.Range("E3").Value = "  -2.81%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'547.66"
$ws.Range("E5").Value = "  -3.96%  "
$ws.Range("D6").Value = "'129.47"
$ws.Range("E6").Value = "  +3.42%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'0.511"
$ws.Range("E8").Value = "  +1.64%  "
$ws.Range("D9").Value = "2.910.03"
$ws.Range("E9").Value = "  -2.76%  "
$ws.Range("E10").Value = "  -5.31%  "
$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").Value = "'0.444"
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").Value = "'4.71"
$ws.Range("E12").Value = "  -6.80%  "
$ws.Range("D13").Value = "'0.0000216"
$ws.Range("E13").Value = "  -2.13%  "
$ws.Range("D14").Value = "'32.57"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("D16").Value = "3.403.02"
$ws.Range("E16").Value = "  -2.78%  "
$ws.Range("D17").Value = "'6.87"
$ws.Range("E17").Value = "  +5.58%  "
$ws.Range("D18").Value = "2.917.47"
$ws.Range("E18").Value = "  -2.79%  "
$ws.Range("D19").Value = "57.517.71"
$ws.Range("E19").Value = "  -4.48%  "
$ws.Range("D20").Value = "'414.91"
$ws.Range("E20").Value = "  -3.58%  "
$ws.Range("D21").Value = "'13.20"
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("D22").Value = "'0.685"
$ws.Range("E22").Value = "  +2.01%  "
$ws.Range("E23").Value = "  -1.60%  "
$ws.Range("D24").Value = "'13.09"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("D25").Value = "'79.44"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  -3.10%  "
$ws.Range("E29").Value = "  +1.28%  "
$ws.Range("E30").Value = "  +2.26%  "
$ws.Range("D31").Value = "'25.15"
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("E32").Value = "  -3.11%  "
$ws.Range("D33").Value = "'0.0965"
$ws.Range("E33").Value = "  +1.49%  "
$ws.Range("E34").Value = "  +0.73%  "
$ws.Range("D35").Value = "'0.927"
$ws.Range("E35").Value = "  -0.70%  "
$ws.Range("E36").Value = "  +0.69%  "
$ws.Range("D37").Value = "'48.11"
$ws.Range("E37").Value = "  -4.34%  "
$ws.Range("D38").Value = "'8.67"
$ws.Range("E38").Value = "  +1.84%  "
$ws.Range("E39").Value = "  +2.53%  "
$ws.Range("E40").Value = "  -0.96%  "
$ws.Range("E41").Value = "  +3.53%  "
$ws.Range("D42").Value = "'374.20"
$ws.Range("E42").Value = "  +0.97%  "
$ws.Range("E43").Value = "  -4.14%  "
$ws.Range("D44").Value = "2.689.55"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "'123.07"
$ws.Range("E46").Value = "  +1.25%  "
$ws.Range("E47").Value = "  +0.73%  "
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("D49").Value = "'1.94"
$ws.Range("E49").Value = "  -2.04%  "
$ws.Range("D50").Value = "'22.87"
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("E51").Value = "  -0.50%  "
